# Update cryptos list - GitHub Actions scheduled refresh
# Price (col D) and Volume(1h) (col E) cells are stored as plain text in the
# source sheet (many prices use a European "thousands dot" style, e.g.
# "26.709.08", so the whole column is text, not numeric). When a new price
# string happens to parse as a normal number (e.g. "18.96"), prefix it with
# a leading apostrophe so Excel keeps storing it as text instead of
# silently converting the cell to a Number, exactly like the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'26.669.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.36%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.630.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.12%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'217.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -1.63%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.13%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -1.51%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -1.35%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'18.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.36%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.05%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "'1.861.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.91%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "'1.627.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.46%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -2.29%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.29%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'63.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.28%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'26.660.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.42%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.0₃0718"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.47%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.09%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'209.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.94%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'4.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.50%  "

# Row 22 - was Toncoin, now Chainlink
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'6.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.96%  "

# Row 23 - was Chainlink, now Toncoin
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "'2.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.99%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  -3.29%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'147.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.07%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -2.39%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  -2.60%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "'15.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.23%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "'0.0498"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.00%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.01%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'3.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.83%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -2.59%  "

# Row 34 - Maker
$ws.Range("D34").Value = "'1.258.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.56%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +0.13%  "

# Row 36 - LidoDAOToken
$ws.Range("D36").Value = "'1.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.64%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  -3.36%  "

# Row 38 - ImmutableX
$ws.Range("D38").Value = "'0.519"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.15%  "

# Row 39 - PaxDollar
$ws.Range("E39").Value = "  +0.06%  "

# Row 40 - ARBITRUM
$ws.Range("D40").Value = "'0.796"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.13%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "'0.797"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.99%  "

# Row 42 - MXToken
$ws.Range("D42").Value = "'2.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.55%  "

# Row 43 - RocketPoolETH
$ws.Range("D43").Value = "'1.775.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.84%  "

# Row 44 - FraxShare
$ws.Range("E44").Value = "  -3.78%  "

# Row 45 - Quant
$ws.Range("D45").Value = "'90.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'59.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.46%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  -2.58%  "

# Row 48 - was BabyDogeCoin, now Cronos
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0517"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.35%  "

# Row 49 - was Cronos, now USDD
$ws.Range("B49").Value = "USDD"
$ws.Range("C49").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D49").Value = "'1.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.16%  "

# Row 50 - was USDD, now Mantle
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.406"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.62%  "

# Row 51 - was Mantle, now EnergySwap
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.30%  "
